# Updates the cryptocurrency price/volume table on the active worksheet
# to match the latest scraped values (GitHub Actions data refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of cell updates: each entry is the target cell reference and its new text value.
$updates = @(
    @{ Ref = 'D2'; Value = '95.688.25' },
    @{ Ref = 'E2'; Value = '  -0.73%  ' },
    @{ Ref = 'D3'; Value = '3.447.88' },
    @{ Ref = 'E3'; Value = '  +4.07%  ' },
    @{ Ref = 'D4'; Value = '1.00' },
    @{ Ref = 'E4'; Value = '  +0.14%  ' },
    @{ Ref = 'D5'; Value = '241.25' },
    @{ Ref = 'E5'; Value = '  -1.91%  ' },
    @{ Ref = 'D6'; Value = '642.95' },
    @{ Ref = 'E6'; Value = '  -0.96%  ' },
    @{ Ref = 'E7'; Value = '  +5.93%  ' },
    @{ Ref = 'D8'; Value = '0.409' },
    @{ Ref = 'E8'; Value = '  +0.15%  ' },
    @{ Ref = 'D9'; Value = '1.00' },
    @{ Ref = 'E9'; Value = '  +0.13%  ' },
    @{ Ref = 'D10'; Value = '0.996' },
    @{ Ref = 'E10'; Value = '  +3.82%  ' },
    @{ Ref = 'D11'; Value = '3.446.31' },
    @{ Ref = 'E11'; Value = '  +4.10%  ' },
    @{ Ref = 'D12'; Value = '43.41' },
    @{ Ref = 'E12'; Value = '  +10.57%  ' },
    @{ Ref = 'D13'; Value = '0.198' },
    @{ Ref = 'E13'; Value = '  -1.97%  ' },
    @{ Ref = 'B14'; Value = 'WrappedBTC' },
    @{ Ref = 'C14'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' },
    @{ Ref = 'D14'; Value = '95.676.08' },
    @{ Ref = 'E14'; Value = '  -0.53%  ' },
    @{ Ref = 'B15'; Value = 'Toncoin' },
    @{ Ref = 'C15'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' },
    @{ Ref = 'D15'; Value = '6.07' },
    @{ Ref = 'E15'; Value = '  +2.77%  ' },
    @{ Ref = 'D16'; Value = '4.096.66' },
    @{ Ref = 'E16'; Value = '  +4.62%  ' },
    @{ Ref = 'E17'; Value = '  +1.90%  ' },
    @{ Ref = 'D18'; Value = '8.57' },
    @{ Ref = 'E18'; Value = '  +1.95%  ' },
    @{ Ref = 'D19'; Value = '3.463.88' },
    @{ Ref = 'E19'; Value = '  +5.27%  ' },
    @{ Ref = 'D20'; Value = '18.23' },
    @{ Ref = 'E20'; Value = '  +10.16%  ' },
    @{ Ref = 'D21'; Value = '11.89' },
    @{ Ref = 'E21'; Value = '  +15.58%  ' },
    @{ Ref = 'D22'; Value = '0.494' },
    @{ Ref = 'E22'; Value = '  +7.97%  ' },
    @{ Ref = 'D23'; Value = '510.43' },
    @{ Ref = 'E23'; Value = '  +4.00%  ' },
    @{ Ref = 'D24'; Value = '3.25' },
    @{ Ref = 'E24'; Value = '  -0.95%  ' },
    @{ Ref = 'B25'; Value = 'NEARProtocol' },
    @{ Ref = 'C25'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' },
    @{ Ref = 'D25'; Value = '6.66' },
    @{ Ref = 'E25'; Value = '  +6.73%  ' },
    @{ Ref = 'B26'; Value = 'PEPE' },
    @{ Ref = 'C26'; Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe' },
    @{ Ref = 'D26'; Value = '0.0000192' },
    @{ Ref = 'E26'; Value = '  -1.42%  ' },
    @{ Ref = 'D27'; Value = '91.97' },
    @{ Ref = 'E27'; Value = '  +0.55%  ' },
    @{ Ref = 'D28'; Value = '12.26' },
    @{ Ref = 'E28'; Value = '  +3.77%  ' },
    @{ Ref = 'D29'; Value = '3.631.97' },
    @{ Ref = 'E29'; Value = '  +4.26%  ' },
    @{ Ref = 'D30'; Value = '11.89' },
    @{ Ref = 'E30'; Value = '  +12.41%  ' },
    @{ Ref = 'E31'; Value = '  -0.29%  ' },
    @{ Ref = 'D32'; Value = '2.76' },
    @{ Ref = 'E32'; Value = '  +14.32%  ' },
    @{ Ref = 'D33'; Value = '0.138' },
    @{ Ref = 'E33'; Value = '  -0.27%  ' },
    @{ Ref = 'E34'; Value = '  +1.24%  ' },
    @{ Ref = 'D35'; Value = '0.581' },
    @{ Ref = 'E35'; Value = '  +9.07%  ' },
    @{ Ref = 'D36'; Value = '30.67' },
    @{ Ref = 'E36'; Value = '  +11.39%  ' },
    @{ Ref = 'E37'; Value = '  +0.28%  ' },
    @{ Ref = 'D38'; Value = '7.78' },
    @{ Ref = 'E38'; Value = '  +5.10%  ' },
    @{ Ref = 'D39'; Value = '1.46' },
    @{ Ref = 'E39'; Value = '  +0.60%  ' },
    @{ Ref = 'D40'; Value = '0.151' },
    @{ Ref = 'E40'; Value = '  +2.67%  ' },
    @{ Ref = 'E41'; Value = '  +0.03%  ' },
    @{ Ref = 'D42'; Value = '0.917' },
    @{ Ref = 'E42'; Value = '  +13.37%  ' },
    @{ Ref = 'D43'; Value = '506.89' },
    @{ Ref = 'E43'; Value = '  +2.30%  ' },
    @{ Ref = 'D44'; Value = '24.19' },
    @{ Ref = 'E44'; Value = '  -1.14%  ' },
    @{ Ref = 'B45'; Value = 'ImmutableX' },
    @{ Ref = 'C45'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' },
    @{ Ref = 'D45'; Value = '1.71' },
    @{ Ref = 'E45'; Value = '  +7.97%  ' },
    @{ Ref = 'B46'; Value = 'VeChain' },
    @{ Ref = 'C46'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Ref = 'D46'; Value = '0.0420' },
    @{ Ref = 'E46'; Value = '  +5.65%  ' },
    @{ Ref = 'D47'; Value = '3.64' },
    @{ Ref = 'E47'; Value = '  -0.34%  ' },
    @{ Ref = 'D48'; Value = '5.54' },
    @{ Ref = 'E48'; Value = '  +4.64%  ' },
    @{ Ref = 'D49'; Value = '3.27' },
    @{ Ref = 'E49'; Value = '  +5.56%  ' },
    @{ Ref = 'D50'; Value = '2.16' },
    @{ Ref = 'E50'; Value = '  +11.80%  ' },
    @{ Ref = 'D51'; Value = '8.25' },
    @{ Ref = 'E51'; Value = '  +0.63%  ' }
)

# Regex that matches strings Excel/IronCalc would otherwise silently coerce into a
# number (plain integers/decimals, optionally signed, optionally in scientific
# notation). Values that match need a leading apostrophe so they are stored as
# literal text (e.g. "1.00" must stay "1.00", not become the number 1).
$numericPattern = '^\s*[+-]?(\d+\.?\d*|\.\d+)([eE][+-]?\d+)?\s*$'

foreach ($update in $updates) {
    $ref = $update.Ref
    $value = $update.Value

    if ($value -match $numericPattern) {
        $value = "'" + $value
    }

    $ws.Range($ref).Value = $value
}
